$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.439.66'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.41'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6950'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07679'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3064'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("E9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.57'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("E10").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.852.06'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("E14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6913'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.321'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.439.72'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008300'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.098.97'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.72'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E21").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.620'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E23").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1492'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("E25").ClearFormats()
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("B26").ClearFormats()
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C26").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.84'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("E26").ClearFormats()
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("B27").ClearFormats()
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C27").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.883'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.93%  '
$ws.Range("E27").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.532'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.244'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.153'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.210'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05111'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7693'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.882'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.148'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.682'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.331.14'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.15%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01872'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("E39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.719'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("E40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9541'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("E41").ClearFormats()
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Quant'
$ws.Range("B42").ClearFormats()
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C42").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.81'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("E42").ClearFormats()
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("B43").ClearFormats()
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C43").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.842'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.50%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.836'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("E45").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.000.27'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5218'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.782'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.32'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.970'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.79%  '
$ws.Range("E51").ClearFormats()
